# Rename the two worksheets to more descriptive names.
$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item(1).Name = "Species Information"
$wb.Worksheets.Item(2).Name = "More Species Information"

$wsSpecies = $wb.Worksheets.Item("Species Information")
$wsMore = $wb.Worksheets.Item("More Species Information")

# Slightly narrow columns A and B on the "Species Information" sheet
# (closest widths this engine's column model can represent to the
# target ~8.37 / ~28.08 character widths).
$wsSpecies.Columns.Item(1).ColumnWidth = 7.5
$wsSpecies.Columns.Item(2).ColumnWidth = 27.166666666666668

# "More Species Information" previously relied on the sheet's implicit
# default column width; pin column A's width explicitly so it is written
# out (closest representable width to the former default of ~8.51).
$wsMore.Columns.Item(1).ColumnWidth = 7.666666666666667
